# Updated ITA_grids model - 2025-08-21 11:48
$wb = $excel.ActiveWorkbook

# --- Sheet "ev_charging_uc": reorder the day/night timeslice lists (C13/C14) ---
$wsEv = $wb.Worksheets.Item("ev_charging_uc")
$wsEv.Range("C13").Value = "RaD,WaP,SaD,FaP,SaP,WaD,FaD,RaP"
$wsEv.Range("C14").Value = "RaN,FaP,SaP,SaN,WaN,FaN,RaP,WaP"

# --- Sheet "re_profiles": rotate the season/ncap_afs lookup table (M4:N7) ---
$wsRe = $wb.Worksheets.Item("re_profiles")
$wsRe.Range("M4").Value = "R"
$wsRe.Range("N4").Value = 0.30301943544655252
$wsRe.Range("M5").Value = "W"
$wsRe.Range("N5").Value = 0.22555529847292916
$wsRe.Range("M6").Value = "S"
$wsRe.Range("N6").Value = 0.40439611291068944
$wsRe.Range("M7").Value = "F"
$wsRe.Range("N7").Value = 0.26702915316982878
